$wb = $excel.ActiveWorkbook

# Insert the new "InvestmentBlock" sheet right before "Lines" (so it lands in the
# same slot the "Lines" sheet used to occupy, pushing Lines/Links one slot right).
$linesSheet = $wb.Worksheets.Item("Lines")
$ws = $wb.Worksheets.Add($linesSheet)
$ws.Name = "InvestmentBlock"

# Match the existing sheets' cell style (centered horizontal alignment, style index 1)
$ws.Range("A1:F5").HorizontalAlignment = -4108

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Meaning"
$ws.Range("C1").Value = "Variable or parameter"
$ws.Range("D1").Value = "Type"
$ws.Range("E1").Value = "Size"
$ws.Range("F1").Value = "Optional"

# Cost
$ws.Range("A2").Value = "Cost"
$ws.Range("B2").Value = "Capital cost"
$ws.Range("C2").Value = "Parameter"
$ws.Range("D2").Value = "float"
$ws.Range("E2").Value = 1

# LowerBound
$ws.Range("A3").Value = "LowerBound"
$ws.Range("B3").Value = "Minimum acceptable size"
$ws.Range("C3").Value = "Parameter"
$ws.Range("D3").Value = "float"
$ws.Range("E3").Value = 1

# UpperBound
$ws.Range("A4").Value = "UpperBound"
$ws.Range("B4").Value = "Maximum acceptable size"
$ws.Range("C4").Value = "Parameter"
$ws.Range("D4").Value = "float"
$ws.Range("E4").Value = 1

# InstalledCapacity
$ws.Range("A5").Value = "InstalledCapacity"
$ws.Range("B5").Value = "Installed capacity"
$ws.Range("C5").Value = "Parameter"
$ws.Range("D5").Value = "float"
$ws.Range("E5").Value = 1

# Column F holds the literal text "False" (not the boolean) on every data row, like
# the rest of the workbook. Copy it from an existing "False" text cell elsewhere so it
# lands as shared-string text (matching style) instead of getting auto-typed as a bool.
$srcFalse = $wb.Worksheets.Item("ThermalUnitBlock").Range("G7")
$srcFalse.Copy($ws.Range("F2"))
$srcFalse.Copy($ws.Range("F3"))
$srcFalse.Copy($ws.Range("F4"))
$srcFalse.Copy($ws.Range("F5"))

# Column widths for the new sheet
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 22.21875
$ws.Columns.Item(3).ColumnWidth = 19.109375
$ws.Columns.Item(4).ColumnWidth = 4.88671875
$ws.Columns.Item(5).ColumnWidth = 4.109375
$ws.Columns.Item(6).ColumnWidth = 8

# Select the full columns A:F (becomes sqref A1:F1048576) and activate this sheet -
# it ends up the active tab, like in the target workbook.
$ws.Range("A:F").Select()
$ws.Activate()

# The BatteryUnitBlock sheet's selection moved to A1:F5 as part of this edit.
$battery = $wb.Worksheets.Item("BatteryUnitBlock")
$battery.Range("A1:F5").Select()

# Re-activate InvestmentBlock so it's the final active sheet/tab.
$ws.Activate()
